$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New members data for "phân công 7" (7th assignment round).
# Entry order follows how the original author actually typed the cells:
# rows 4-6 filled MSSV+name together, but rows 7-8 had both MSSV values
# (C7, C8) entered before either name (D7, D8).
$ws.Cells.Item(4, 3).Value = "0712187"
$ws.Cells.Item(4, 4).Value = "Lý Hoài"

$ws.Cells.Item(5, 3).Value = "0712188"
$ws.Cells.Item(5, 4).Value = "Phan Lê Huỳnh"

$ws.Cells.Item(6, 3).Value = "0712236"
$ws.Cells.Item(6, 4).Value = "Phan Vũ Lâm"

$ws.Cells.Item(7, 3).Value = "0712365"
$ws.Cells.Item(8, 3).Value = "0712381"
$ws.Cells.Item(7, 4).Value = "Nguyễn Hồ Mẫn Sáng"
$ws.Cells.Item(8, 4).Value = "La Minh Tâm"

for ($r = 4; $r -le 8; $r++) {
    for ($col = 5; $col -le 11; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.Value = 1
        # Match the "Neutral" cell-style look (Times New Roman 12, brownish
        # font colour, light-yellow fill, thin border) already used for
        # columns C/D in this row, plus a 0% number format.
        $cell.Font.Name = "Times New Roman"
        $cell.Font.Size = 12
        $cell.Font.Color = 26012
        $cell.Interior.Color = 10284031
        $cell.Borders.LineStyle = 1
        $cell.NumberFormat = "0%"
    }
}

# Adjust column D width to match new content (no longer best-fit)
$ws.Columns.Item(4).ColumnWidth = 20.6

# Update the active selection
$ws.Range("H10").Select()
